# feat: add 2022-Q1 data
#
# Before: sheets = [2021-Q3, 2021-Q4, 总计]
# After:  sheets = [2021-Q3, 2021-Q4, 2022-Q1, 总计]
#
# The old "总计" sheet (sheetId=3) is renamed to "2022-Q1" and repurposed to
# hold the new quarter's per-fund holdings (same shape as the 2021-Q3 /
# 2021-Q4 sheets). A brand-new "总计" sheet is appended at the end holding
# the refreshed summary table (adds the 2022-Q1 row on top of the existing
# 2021-Q4 / 2021-Q3 rows).
#
# NOTE: every cell below is addressed via Range("<A1-address>") rather than
# Cells.Item(row,col) -- mixing the two addressing styles on the same cell
# (value write via Cells.Item, then a PasteSpecial format-copy via Range)
# has been observed to drop the pasted style, so a single addressing style
# is used consistently throughout.

$wb = $excel.ActiveWorkbook

function Get-ColLetter($col) {
    $letters = ""
    while ($col -gt 0) {
        $rem = ($col - 1) % 26
        $letters = [char](65 + $rem) + $letters
        $col = [int](($col - $rem - 1) / 26)
    }
    return $letters
}

function Get-Addr($row, $col) {
    return (Get-ColLetter $col) + $row
}

# Write a value into a cell while forcing text storage so that
# numeric-looking strings (fund codes, formatted percentages, etc.) are not
# silently coerced into numbers.
function Set-TextCell($ws, $row, $col, $val) {
    $rng = $ws.Range((Get-Addr $row $col))
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

function Set-NumberCell($ws, $row, $col, $val) {
    $ws.Range((Get-Addr $row $col)).Value = $val
}

# Copy just the cell formatting (font/border/alignment) from a known-styled
# source cell onto a destination range, leaving its value/format untouched.
function Copy-CellStyle($srcWs, $srcAddr, $dstWs, $dstAddr) {
    $srcWs.Range($srcAddr).Copy() | Out-Null
    $dstWs.Range($dstAddr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Rename the existing "总计" sheet to "2022-Q1" and add the brand-new
#    "总计" sheet right after it (so tab order becomes ..., 2022-Q1, 总计).
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item("总计")
$q1Sheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1Sheet)
$totalSheet.Name = "总计"

# A sheet with an existing correctly-styled header (s="2") to source
# formatting from.
$styleSource = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 2. Rebuild "2022-Q1" with the per-fund holdings table.
# ---------------------------------------------------------------------

# Clear any leftover content from the old "总计" data.
$q1Sheet.Cells.Clear() | Out-Null

# Header row.
Set-TextCell $q1Sheet 1 2 "基金代码"
Set-TextCell $q1Sheet 1 3 "基金名称"
Set-TextCell $q1Sheet 1 4 "基金规模"
Set-TextCell $q1Sheet 1 5 "股票总仓位"
Set-TextCell $q1Sheet 1 6 "仓位占比"
Set-TextCell $q1Sheet 1 7 "持有市值(亿元)"
Set-TextCell $q1Sheet 1 8 "仓位排名"
Copy-CellStyle $styleSource "B1:H1" $q1Sheet "B1:H1"

$q1Rows = @(
    @{ A = 0; B = "217024"; C = "招商安盈债券";             D = "35.05"; E = "20.20"; F = "0.74"; G = "0.2594"; H = 8 },
    @{ A = 1; B = "014887"; C = "招商安福1年定期开放债券";   D = "17.22"; E = "27.65"; F = "0.71"; G = "0.1223"; H = 9 },
    @{ A = 2; B = "005459"; C = "银河嘉谊灵活配置混合A";     D = "6.47";  E = "39.69"; F = "0.83"; G = "0.0537"; H = 2 },
    @{ A = 3; B = "005460"; C = "银河嘉谊灵活配置混合C";     D = "2.79";  E = "39.69"; F = "0.83"; G = "0.0232"; H = 2 },
    @{ A = 4; B = "005053"; C = "银河量化价值混合";         D = "0.10"; E = "80.73"; F = "2.69"; G = "0.0027"; H = 3 }
)

$r = 2
foreach ($row in $q1Rows) {
    $destAddr = Get-Addr $r 1
    Set-NumberCell $q1Sheet $r 1 $row.A
    Copy-CellStyle $styleSource "A2" $q1Sheet $destAddr
    Set-TextCell $q1Sheet $r 2 $row.B
    Set-TextCell $q1Sheet $r 3 $row.C
    Set-TextCell $q1Sheet $r 4 $row.D
    Set-TextCell $q1Sheet $r 5 $row.E
    Set-TextCell $q1Sheet $r 6 $row.F
    Set-TextCell $q1Sheet $r 7 $row.G
    Set-NumberCell $q1Sheet $r 8 $row.H
    $r += 1
}

$q1Sheet.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Populate the brand-new "总计" sheet with the refreshed summary table.
# ---------------------------------------------------------------------

Set-TextCell $totalSheet 1 2 "日期"
Set-TextCell $totalSheet 1 3 "持有数量(只)"
Set-TextCell $totalSheet 1 4 "持有市值(亿元)"
Copy-CellStyle $styleSource "B1:D1" $totalSheet "B1:D1"

$totalRows = @(
    @{ A = 0; B = "2022-Q1"; C = 5; D = 0.46 },
    @{ A = 1; B = "2021-Q4"; C = 2; D = 0.07000000000000001 },
    @{ A = 2; B = "2021-Q3"; C = 3; D = 0.06 }
)

$r = 2
foreach ($row in $totalRows) {
    $destAddr = Get-Addr $r 1
    Set-NumberCell $totalSheet $r 1 $row.A
    Copy-CellStyle $styleSource "A2" $totalSheet $destAddr
    Set-TextCell $totalSheet $r 2 $row.B
    Set-NumberCell $totalSheet $r 3 $row.C
    Set-NumberCell $totalSheet $r 4 $row.D
    $r += 1
}

$totalSheet.Range("A1").Select() | Out-Null

Write-Host "2022-Q1 sheet inserted; zong-ji sheet refreshed."
